$d = $word.ActiveDocument

# =========================================================================
# Edit region 1 (paragraph containing "...pose l<m>or</m> en feuille &...")
#   "or"            -> "or en feuille"     (run B: plain Arial text run)
#   "</m>" run      -> drop explicit Bold=False (w:b val="0" is removed)
#   " en feuille "  -> " "                 (run D: plain Arial text run)
# =========================================================================

$anchor1 = $d.Content
$found1 = $anchor1.Find.Execute("pose l<m>or</m> en feuille ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Edit region 1 anchor text not found"
}
$base1 = $anchor1.Start

# "or" run -> "or en feuille"
$rOr = $d.Range($base1 + 9, $base1 + 11)
if ($rOr.Text -ne "or") { throw "Edit region 1: unexpected text [$($rOr.Text)] where 'or' was expected" }
$rOr.Text = "or en feuille"

# "</m>" run -> drop Bold (the text itself is unchanged)
$rCloseM = $d.Range($rOr.End, $rOr.End + 4)
if ($rCloseM.Text -ne "</m>") { throw "Edit region 1: unexpected text [$($rCloseM.Text)] where '</m>' was expected" }
$rCloseM.Font.Bold = $false

# " en feuille " run -> " "
$rEnFeuille = $d.Range($rCloseM.End, $rCloseM.End + 12)
if ($rEnFeuille.Text -ne " en feuille ") { throw "Edit region 1: unexpected text [$($rEnFeuille.Text)] where ' en feuille ' was expected" }
$rEnFeuille.Text = " "

# =========================================================================
# Edit region 2 (paragraph containing "...en <env>une muraille</env></head>")
#   "en "           -> "en une "           (run B: plain Arial text run)
#   "<env>" run     -> drop explicit Bold=False (w:b val="0" is removed)
#   "une muraille"  -> "muraille"          (run D: plain Arial text run)
# =========================================================================

$anchor2 = $d.Content
$found2 = $anchor2.Find.Execute("en <env>une muraille</env>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Edit region 2 anchor text not found"
}
$base2 = $anchor2.Start

# "en " run -> "en une "
$rEn = $d.Range($base2, $base2 + 3)
if ($rEn.Text -ne "en ") { throw "Edit region 2: unexpected text [$($rEn.Text)] where 'en ' was expected" }
$rEn.Text = "en une "

# "<env>" run -> drop Bold (the text itself is unchanged)
$rOpenEnv = $d.Range($rEn.End, $rEn.End + 5)
if ($rOpenEnv.Text -ne "<env>") { throw "Edit region 2: unexpected text [$($rOpenEnv.Text)] where '<env>' was expected" }
$rOpenEnv.Font.Bold = $false

# "une muraille" run -> "muraille"
$rUneMuraille = $d.Range($rOpenEnv.End, $rOpenEnv.End + 12)
if ($rUneMuraille.Text -ne "une muraille") { throw "Edit region 2: unexpected text [$($rUneMuraille.Text)] where 'une muraille' was expected" }
$rUneMuraille.Text = "muraille"

Write-Output "done"
